# får ikke til sortering av dagsverk
# Append run/trip numbers to the shift-time (and XX/OO/TT) cells on both
# week sheets, and widen the "XX / OO / TT" conditional-format match so it
# still fires now that those short codes carry a trailing space.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Cell value updates — OSL_01
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("OSL_01")

$ws1.Range("D2").Value2 = "OO 1429"
$ws1.Range("E2").Value2 = "9:40 - 19:06 5005-73"
$ws1.Range("F2").Value2 = "14:01 - 21:42 1235"
$ws1.Range("G2").Value2 = "16:07 - 22:38 9306"
$ws1.Range("H2").Value2 = "13:00 - 21:00 "

$ws1.Range("D3").Value2 = "6:00 - 14:00 905301500"
$ws1.Range("E3").Value2 = "5:30 - 15:00 9301"
$ws1.Range("F3").Value2 = "5:00 - 13:00 "
$ws1.Range("H3").Value2 = "OO "

$ws1.Range("D4").Value2 = "15:00 - 23:00 9308"
$ws1.Range("E4").Value2 = "14:00 - 22:00 9329-X"
$ws1.Range("F4").Value2 = "15:00 - 23:00 "
$ws1.Range("G4").Value2 = "23:00 - 7:00 9312"
$ws1.Range("H4").Value2 = "23:45 - 7:36 1775"

$ws1.Range("D5").Value2 = "8:01 - 17:29 5014"
$ws1.Range("E5").Value2 = "7:13 - 15:01 1424-Mod1"
$ws1.Range("F5").Value2 = "7:17 - 15:17 "
$ws1.Range("H5").Value2 = "OO "

$ws1.Range("D6").Value2 = "21:39 - 5:45 1567"
$ws1.Range("F6").Value2 = "OO 1607"
$ws1.Range("G6").Value2 = "7:35 - 15:47 1705"
$ws1.Range("H6").Value2 = "7:43 - 16:20 3008"

$ws1.Range("D7").Value2 = "14:00 - 22:00 "
$ws1.Range("E7").Value2 = "TT 9906001500-H"
$ws1.Range("F7").Value2 = "6:00 - 15:00 "
$ws1.Range("H7").Value2 = "OO "

# ---------------------------------------------------------------------
# 2) Cell value updates — OSL_02
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("OSL_02")

$ws2.Range("B2").Value2 = "15:00 - 23:00 9309"
$ws2.Range("C2").Value2 = "23:23 - 6:28 3119"
$ws2.Range("F2").Value2 = "OO 9323-X"
$ws2.Range("G2").Value2 = "6:00 - 14:00 3031"
$ws2.Range("H2").Value2 = "13:28 - 22:38 5004"

$ws2.Range("B3").Value2 = "15:08 - 22:26 5005-73"
$ws2.Range("C3").Value2 = "14:01 - 21:42 1447"
$ws2.Range("D3").Value2 = "14:45 - 22:56 9330-X"
$ws2.Range("E3").Value2 = "15:00 - 23:00 "
$ws2.Range("F3").Value2 = "TT "
$ws2.Range("H3").Value2 = "OO "

$ws2.Range("B4").Value2 = "OO 906001600"
$ws2.Range("C4").Value2 = "6:00 - 16:00 9322-X"
$ws2.Range("D4").Value2 = "5:00 - 13:00 "
$ws2.Range("E4").Value2 = "XX 9323-X"
$ws2.Range("F4").Value2 = "6:00 - 14:00 9305"
$ws2.Range("G4").Value2 = "7:00 - 15:00 3134"
$ws2.Range("H4").Value2 = "10:35 - 19:38 3006"

$ws2.Range("B5").Value2 = "13:31 - 21:23 "
$ws2.Range("C5").Value2 = "21:23 - 7:28 3109"
$ws2.Range("D5").Value2 = "22:00 - 6:00 9331-X"
$ws2.Range("F5").Value2 = "7:00 - 16:00 "
$ws2.Range("H5").Value2 = "OO "

$ws2.Range("C6").Value2 = "OO 906001600"
$ws2.Range("D6").Value2 = "6:00 - 16:00 5002"
$ws2.Range("E6").Value2 = "6:33 - 13:13 1458-Mod1"
$ws2.Range("F6").Value2 = "15:12 - 22:20 3020"
$ws2.Range("G6").Value2 = "9:41 - 20:08 3132"
$ws2.Range("H6").Value2 = "15:18 - 22:28 "

$ws2.Range("B7").Value2 = "23:00 - 7:00 9312"
$ws2.Range("C7").Value2 = "23:00 - 7:00 9312"
$ws2.Range("E7").Value2 = "TT 3101"
$ws2.Range("F7").Value2 = "7:21 - 14:33 "
$ws2.Range("H7").Value2 = "OO "

# ---------------------------------------------------------------------
# 3) Conditional formatting: the "is it XX / OO / TT" rule on every
#    B2:H7 cell (both sheets) must match with a trailing space now,
#    since those codes are written as "XX ", "OO ", "TT ".
#    On every one of these cells that rule is the LAST FormatCondition.
# ---------------------------------------------------------------------
$cols = @("B", "C", "D", "E", "F", "G", "H")
$rows = @(2, 3, 4, 5, 6, 7)

foreach ($ws in @($ws1, $ws2)) {
    foreach ($col in $cols) {
        foreach ($row in $rows) {
            $ref = "$col$row"
            $fcs = $ws.Range($ref).FormatConditions
            $fc = $fcs.Item($fcs.Count)
            $fc.Formula1 = '=(' + $ref + '="XX ")OR (' + $ref + '="OO ")OR (' + $ref + '="TT ")'
        }
    }
}
